# Postgame hitter report update - "getting closer on postgame hitter"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pitch block 1 (rows 9-17) ---
$ws.Range("J10").Value = 2
$ws.Range("M10").Value = ""

$ws.Range("M12").Value = ""

$ws.Range("J14").Value = "Roblez"
$ws.Range("M14").Value = "Ground Ball"

$ws.Range("M15").Value = "Out"

$ws.Range("J16").Value = "88-90 MPH"

$ws.Range("J17").Value = "FB,CB,CH"

# --- Pitch block 2 (rows 18-26) ---
$ws.Range("J19").Value = 3
$ws.Range("M19").Value = ""

$ws.Range("J20").Value = 2

$ws.Range("M21").Value = ""

$ws.Range("M23").Value = "Popup"

$ws.Range("J26").Value = "FB,CB,CH"

# --- Pitch block 3 (rows 27-35) ---
$ws.Range("M28").Value = ""

$ws.Range("M30").Value = ""

$ws.Range("J35").Value = "SL,FB,CB,CH"

# --- Pitch block 4 (rows 36-44) ---
$ws.Range("J37").Value = 6
$ws.Range("M37").Value = ""

$ws.Range("J38").Value = 1

$ws.Range("M39").Value = ""

$ws.Range("J41").Value = "Herbst"
$ws.Range("M41").Value = "Line Drive"

$ws.Range("M42").Value = "Single"

$ws.Range("J43").Value = "83-85 MPH"

$ws.Range("J44").Value = "SL,FB,CB,CH"
